$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Edit 1: shape 4 ("Thank you for making a new cheatsheet...") ---
# Paragraph 5: "If the execution is canceled, then the strategy will throw TimeoutRejectedException."
# becomes: "If the execution is cancelled, then the strategy will throw TimeoutRejectedException."
# The three runs covering "If the execution is " / "canceled" / ", then the strategy will throw "
# get merged into a single run (keeping the first run's formatting).
$shp1 = $s.Shapes.Item(4)
$tr1 = $shp1.TextFrame.TextRange
$para5 = $tr1.Paragraphs(5)

$toDelete = $para5.Characters(21, 39)   # "canceled, then the strategy will throw "
$toDelete.Delete()

$anchor1 = $para5.Characters(1, 20)     # "If the execution is " (still carries run-1 formatting)
$anchor1.InsertAfter("cancelled, then the strategy will throw ") | Out-Null

# --- Edit 2: shape 13 ("Use a layout that flows...", id=24) ---
# "Specify asynchronously delegate for timeout notification "
# becomes: "Specify asynchronous delegate for notification "
$shp2 = $s.Shapes.Item(13)
$tr2 = $shp2.TextFrame.TextRange
$para1 = $tr2.Paragraphs(1)

$prefix = $para1.Characters(1, 1)       # "S" (keeps run formatting)
$rest = $para1.Characters(2, $para1.Length - 1)
$rest.Delete()
$prefix.InsertAfter("pecify asynchronous delegate for notification ") | Out-Null
